$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2 through 387 all change from 2023-09-20 (45189)
# to 2023-09-21 (45190).
$ws.Range("C2:C387").Value = 45190
